$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1242.5  # H17
$ws.Cells.Item(17, 9).Value = 500  # I17
$ws.Cells.Item(17, 10).Value = 1348.5714  # J17
$ws.Cells.Item(17, 11).Value = 1500  # K17
$ws.Cells.Item(17, 12).Value = 4045.7142  # L17
$ws.Cells.Item(17, 13).Value = -1332  # M17
$ws.Cells.Item(17, 14).Value = -4381.7142  # N17
$ws.Cells.Item(70, 8).Value = 885.1707  # H70
$ws.Cells.Item(70, 9).Value = 928.2222  # I70
$ws.Cells.Item(70, 10).Value = 802.1429000000001  # J70
$ws.Cells.Item(70, 11).Value = 2784.6666  # K70
$ws.Cells.Item(70, 12).Value = 2406.4287  # L70
$ws.Cells.Item(70, 13).Value = -2514.6666  # M70
$ws.Cells.Item(70, 14).Value = -2946.4287  # N70
$ws.Cells.Item(73, 8).Value = 885.1707  # H73
$ws.Cells.Item(73, 9).Value = 928.2222  # I73
$ws.Cells.Item(73, 10).Value = 802.1429000000001  # J73
$ws.Cells.Item(73, 11).Value = 2784.6666  # K73
$ws.Cells.Item(73, 12).Value = 2406.4287  # L73
$ws.Cells.Item(73, 13).Value = -1848.6666  # M73
$ws.Cells.Item(73, 14).Value = -4278.4287  # N73
$ws.Cells.Item(134, 8).Value = 111183420  # H134
$ws.Cells.Item(134, 10).Value = 111183420  # J134
$ws.Cells.Item(134, 12).Value = 111183420  # L134
$ws.Cells.Item(134, 14).Value = -111193560  # N134
$ws.Cells.Item(138, 8).Value = 2102.3066  # H138
$ws.Cells.Item(138, 9).Value = 1140.75  # I138
$ws.Cells.Item(138, 10).Value = 3201.2285  # J138
$ws.Cells.Item(138, 11).Value = 3422.25  # K138
$ws.Cells.Item(138, 12).Value = 9603.6855  # L138
$ws.Cells.Item(138, 13).Value = 1717.75  # M138
$ws.Cells.Item(138, 14).Value = -19883.6855  # N138
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10424.517  # H32
$ws.Cells.Item(32, 9).Value = 8067.3657  # I32
$ws.Cells.Item(32, 10).Value = 27996  # J32
$ws.Cells.Item(32, 11).Value = 8067.3657  # K32
$ws.Cells.Item(32, 12).Value = 27996  # L32
$ws.Cells.Item(32, 13).Value = -7780.3657  # M32
$ws.Cells.Item(32, 14).Value = -28570  # N32
$ws.Cells.Item(102, 8).Value = 2589.0908  # H102
$ws.Cells.Item(102, 9).Value = 2171.4285  # I102
$ws.Cells.Item(102, 10).Value = 3320  # J102
$ws.Cells.Item(102, 11).Value = 2171.4285  # K102
$ws.Cells.Item(102, 12).Value = 3320  # L102
$ws.Cells.Item(102, 13).Value = -549.4285  # M102
$ws.Cells.Item(102, 14).Value = -6564  # N102
$ws.Cells.Item(132, 8).Value = 920464.9399999999  # H132
$ws.Cells.Item(132, 9).Value = 1445789.9  # I132
$ws.Cells.Item(132, 10).Value = 8058.421  # J132
$ws.Cells.Item(132, 11).Value = 4337369.699999999  # K132
$ws.Cells.Item(132, 12).Value = 24175.263  # L132
$ws.Cells.Item(132, 13).Value = -4334839.699999999  # M132
$ws.Cells.Item(132, 14).Value = -29235.263  # N132
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1754.931  # H20
$ws.Cells.Item(20, 9).Value = 949.2778  # I20
$ws.Cells.Item(20, 10).Value = 3073.2727  # J20
$ws.Cells.Item(20, 11).Value = 949.2778  # K20
$ws.Cells.Item(20, 12).Value = 3073.2727  # L20
$ws.Cells.Item(20, 13).Value = -702.2778  # M20
$ws.Cells.Item(20, 14).Value = -3567.2727  # N20
$ws.Cells.Item(99, 8).Value = 2249.0908  # H99
$ws.Cells.Item(99, 9).Value = 0  # I99
$ws.Cells.Item(99, 10).Value = 2249.0908  # J99
$ws.Cells.Item(99, 11).Value = 0  # K99
$ws.Cells.Item(99, 12).Value = 2249.0908  # L99
$ws.Cells.Item(99, 13).ClearContents()  # M99 removed (was 498.5)
$ws.Cells.Item(99, 14).Value = -5245.0908  # N99
$ws.Cells.Item(105, 8).Value = 1854.9474  # H105
$ws.Cells.Item(105, 9).Value = 1427  # I105
$ws.Cells.Item(105, 11).Value = 1427  # K105
$ws.Cells.Item(105, 13).Value = 320  # M105
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1771.9692  # H31
$ws.Cells.Item(31, 9).Value = 1087.1132  # I31
$ws.Cells.Item(31, 10).Value = 4796.75  # J31
$ws.Cells.Item(31, 11).Value = 1087.1132  # K31
$ws.Cells.Item(31, 12).Value = 4796.75  # L31
$ws.Cells.Item(31, 13).Value = -792.1132  # M31
$ws.Cells.Item(31, 14).Value = -5386.75  # N31
$ws.Cells.Item(34, 8).Value = 1771.9692  # H34
$ws.Cells.Item(34, 9).Value = 1087.1132  # I34
$ws.Cells.Item(34, 10).Value = 4796.75  # J34
$ws.Cells.Item(34, 11).Value = 1087.1132  # K34
$ws.Cells.Item(34, 12).Value = 4796.75  # L34
$ws.Cells.Item(34, 13).Value = -885.1132  # M34
$ws.Cells.Item(34, 14).Value = -5200.75  # N34
$ws.Cells.Item(43, 8).Value = 17134.572  # H43
$ws.Cells.Item(43, 10).Value = 17134.572  # J43
$ws.Cells.Item(43, 12).Value = 17134.572  # L43
$ws.Cells.Item(43, 14).Value = -17502.572  # N43
$ws.Cells.Item(99, 8).Value = 7145650.5  # H99
$ws.Cells.Item(99, 9).Value = 2709.2  # I99
$ws.Cells.Item(99, 10).Value = 25003004  # J99
$ws.Cells.Item(99, 11).Value = 2709.2  # K99
$ws.Cells.Item(99, 12).Value = 25003004  # L99
$ws.Cells.Item(99, 13).Value = -1211.2  # M99
$ws.Cells.Item(99, 14).Value = -25006000  # N99
$ws.Cells.Item(101, 8).Value = 17134.572  # H101
$ws.Cells.Item(101, 10).Value = 17134.572  # J101
$ws.Cells.Item(101, 12).Value = 17134.572  # L101
$ws.Cells.Item(101, 14).Value = -23624.572  # N101
$ws.Cells.Item(126, 8).Value = 7145650.5  # H126
$ws.Cells.Item(126, 9).Value = 2709.2  # I126
$ws.Cells.Item(126, 10).Value = 25003004  # J126
$ws.Cells.Item(126, 11).Value = 8127.599999999999  # K126
$ws.Cells.Item(126, 12).Value = 75009012  # L126
$ws.Cells.Item(126, 13).Value = -5657.599999999999  # M126
$ws.Cells.Item(126, 14).Value = -75013952  # N126
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1257.0714  # H34
$ws.Cells.Item(34, 10).Value = 1399.9166  # J34
$ws.Cells.Item(34, 12).Value = 4199.7498  # L34
$ws.Cells.Item(34, 14).Value = -4367.7498  # N34
$ws.Cells.Item(132, 8).Value = 838.6667  # H132
$ws.Cells.Item(132, 9).Value = 755.5909  # I132
$ws.Cells.Item(132, 10).Value = 1752.5  # J132
$ws.Cells.Item(132, 11).Value = 6800.3181  # K132
$ws.Cells.Item(132, 12).Value = 15772.5  # L132
$ws.Cells.Item(132, 13).Value = -4270.3181  # M132
$ws.Cells.Item(132, 14).Value = -20832.5  # N132
$ws.Cells.Item(136, 8).Value = 2440  # H136
$ws.Cells.Item(136, 9).Value = 2096.6667  # I136
$ws.Cells.Item(136, 11).Value = 6290.000100000001  # K136
$ws.Cells.Item(136, 13).Value = -1190.000100000001  # M136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10899.8  # H70
$ws.Cells.Item(70, 9).Value = 19550  # I70
$ws.Cells.Item(70, 10).Value = 5133  # J70
$ws.Cells.Item(70, 11).Value = 19550  # K70
$ws.Cells.Item(70, 12).Value = 5133  # L70
$ws.Cells.Item(70, 13).Value = -19280  # M70
$ws.Cells.Item(70, 14).Value = -5673  # N70
$ws.Cells.Item(73, 8).Value = 10899.8  # H73
$ws.Cells.Item(73, 9).Value = 19550  # I73
$ws.Cells.Item(73, 10).Value = 5133  # J73
$ws.Cells.Item(73, 11).Value = 19550  # K73
$ws.Cells.Item(73, 12).Value = 5133  # L73
$ws.Cells.Item(73, 13).Value = -18614  # M73
$ws.Cells.Item(73, 14).Value = -7005  # N73
$ws.Cells.Item(95, 8).Value = 8970.571  # H95
$ws.Cells.Item(95, 10).Value = 8970.571  # J95
$ws.Cells.Item(95, 12).Value = 8970.571  # L95
$ws.Cells.Item(95, 14).Value = -14462.571  # N95
$ws.Cells.Item(99, 8).Value = 14005.462  # H99
$ws.Cells.Item(99, 9).Value = 10638.714  # I99
$ws.Cells.Item(99, 11).Value = 10638.714  # K99
$ws.Cells.Item(99, 13).Value = -8392.714  # M99
$ws.Cells.Item(132, 8).Value = 4268.6  # H132
$ws.Cells.Item(132, 9).Value = 3369  # I132
$ws.Cells.Item(132, 10).Value = 4868.3335  # J132
$ws.Cells.Item(132, 11).Value = 10107  # K132
$ws.Cells.Item(132, 12).Value = 14605.0005  # L132
$ws.Cells.Item(132, 13).Value = -7577  # M132
$ws.Cells.Item(132, 14).Value = -19665.0005  # N132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1858.8572  # H82
$ws.Cells.Item(82, 9).Value = 1432  # I82
$ws.Cells.Item(82, 10).Value = 2121.5386  # J82
$ws.Cells.Item(82, 11).Value = 1432  # K82
$ws.Cells.Item(82, 12).Value = 2121.5386  # L82
$ws.Cells.Item(82, 13).Value = -1071  # M82
$ws.Cells.Item(82, 14).Value = -2843.5386  # N82
$ws.Cells.Item(85, 8).Value = 1858.8572  # H85
$ws.Cells.Item(85, 9).Value = 1432  # I85
$ws.Cells.Item(85, 10).Value = 2121.5386  # J85
$ws.Cells.Item(85, 11).Value = 1432  # K85
$ws.Cells.Item(85, 12).Value = 2121.5386  # L85
$ws.Cells.Item(85, 13).Value = -184  # M85
$ws.Cells.Item(85, 14).Value = -4617.5386  # N85
$ws.Cells.Item(97, 8).Value = 11945.692  # H97
$ws.Cells.Item(97, 10).Value = 11945.692  # J97
$ws.Cells.Item(97, 12).Value = 11945.692  # L97
$ws.Cells.Item(97, 14).Value = -13927.692  # N97
$ws.Cells.Item(132, 8).Value = 37645.465  # H132
$ws.Cells.Item(132, 9).Value = 70433.39999999999  # I132
$ws.Cells.Item(132, 10).Value = 4857.533  # J132
$ws.Cells.Item(132, 11).Value = 211300.2  # K132
$ws.Cells.Item(132, 12).Value = 14572.599  # L132
$ws.Cells.Item(132, 13).Value = -208770.2  # M132
$ws.Cells.Item(132, 14).Value = -19632.599  # N132
$ws.Cells.Item(136, 8).Value = 1544.0312  # H136
$ws.Cells.Item(136, 9).Value = 1353.04  # I136
$ws.Cells.Item(136, 10).Value = 2226.1428  # J136
$ws.Cells.Item(136, 11).Value = 4059.12  # K136
$ws.Cells.Item(136, 12).Value = 6678.428400000001  # L136
$ws.Cells.Item(136, 13).Value = -1509.12  # M136
$ws.Cells.Item(136, 14).Value = -11778.4284  # N136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3822.3333  # H62
$ws.Cells.Item(62, 9).Value = 2862.5  # I62
$ws.Cells.Item(62, 10).Value = 4171.364  # J62
$ws.Cells.Item(62, 11).Value = 2862.5  # K62
$ws.Cells.Item(62, 12).Value = 4171.364  # L62
$ws.Cells.Item(62, 13).Value = -2238.5  # M62
$ws.Cells.Item(62, 14).Value = -5419.364  # N62
$ws.Cells.Item(65, 8).Value = 3822.3333  # H65
$ws.Cells.Item(65, 9).Value = 2862.5  # I65
$ws.Cells.Item(65, 10).Value = 4171.364  # J65
$ws.Cells.Item(65, 11).Value = 14312.5  # K65
$ws.Cells.Item(65, 12).Value = 20856.82  # L65
$ws.Cells.Item(65, 13).Value = -11192.5  # M65
$ws.Cells.Item(65, 14).Value = -27096.82  # N65
$ws.Cells.Item(97, 8).Value = 19535.334  # H97
$ws.Cells.Item(97, 10).Value = 19535.334  # J97
$ws.Cells.Item(97, 12).Value = 19535.334  # L97
$ws.Cells.Item(97, 14).Value = -21517.334  # N97
$ws.Cells.Item(126, 8).Value = 3214.8147  # H126
$ws.Cells.Item(126, 9).Value = 2902.1428  # I126
$ws.Cells.Item(126, 10).Value = 4309.1665  # J126
$ws.Cells.Item(126, 11).Value = 8706.428400000001  # K126
$ws.Cells.Item(126, 12).Value = 12927.4995  # L126
$ws.Cells.Item(126, 13).Value = -6236.428400000001  # M126
